$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting for numeric-looking values
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '96.402.76'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '3.698.03'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '236.46'
$ws.Range('E5').Value = '  -2.56%  '
$ws.Range('D6').Value = '1.90'
$ws.Range('E6').Value = '  +3.41%  '
$ws.Range('D7').Value = '650.82'
$ws.Range('E7').Value = '  -0.79%  '
$ws.Range('D8').Value = '0.427'
$ws.Range('E8').Value = '  +1.33%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('D11').Value = '3.698.04'
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('D12').Value = '44.26'
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('D14').Value = '0.0000291'
$ws.Range('E14').Value = '  +12.65%  '
$ws.Range('D15').Value = '6.72'
$ws.Range('E15').Value = '  +2.89%  '
$ws.Range('D16').Value = '4.386.07'
$ws.Range('E16').Value = '  +1.48%  '
$ws.Range('D17').Value = '97.216.86'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '8.82'
$ws.Range('E18').Value = '  +13.84%  '
$ws.Range('D19').Value = '3.696.22'
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('D20').Value = '12.92'
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('D21').Value = '18.74'
$ws.Range('E21').Value = '  +1.79%  '
$ws.Range('D22').Value = '0.503'
$ws.Range('E22').Value = '  -5.94%  '
$ws.Range('D23').Value = '518.07'
$ws.Range('E23').Value = '  +0.96%  '
$ws.Range('E24').Value = '  -2.66%  '
$ws.Range('D25').Value = '0.0000205'
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('D26').Value = '6.97'
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('D27').Value = '101.00'
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').Value = '13.12'
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').Value = '0.172'
$ws.Range('E29').Value = '  +2.82%  '
$ws.Range('D30').Value = '3.01'
$ws.Range('E30').Value = '  -1.28%  '
$ws.Range('D31').Value = '12.04'
$ws.Range('E31').Value = '  +0.92%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '1.84'
$ws.Range('E33').Value = '  +4.50%  '
$ws.Range('B34').Value = 'Cronos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D34').Value = '0.184'
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').Value = '655.89'
$ws.Range('E36').Value = '  +6.43%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = '32.19'
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('D39').Value = '8.78'
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').Value = '2.09'
$ws.Range('E41').Value = '  +6.99%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '6.85'
$ws.Range('E42').Value = '  +10.99%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '41.02'
$ws.Range('E43').Value = '  -4.51%  '
$ws.Range('E44').Value = '  +0.85%  '
$ws.Range('D45').Value = '0.963'
$ws.Range('E45').Value = '  +0.64%  '
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('D47').Value = '0.432'
$ws.Range('E47').Value = '  +5.23%  '
$ws.Range('D48').Value = '23.57'
$ws.Range('E48').Value = '  -0.09%  '
$ws.Range('E49').Value = '  -1.97%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '8.48'
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('B51').Value = 'MantraDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D51').Value = '3.53'
$ws.Range('E51').Value = '  +2.25%  '
